$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.289.62'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '3.283.67'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'579.60"
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = "'179.22"
$ws.Range('E6').Value = '  -2.89%  '
$ws.Range('D7').Value = "'0.629"
$ws.Range('E7').Value = '  +4.12%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -2.49%  '
$ws.Range('D10').Value = "'6.73"
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('D11').Value = "'0.402"
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').Value = '3.850.61'
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').Value = '66.300.59'
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('D15').Value = "'26.37"
$ws.Range('E15').Value = '  -3.76%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.340.50'
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = "'0.0000164"
$ws.Range('E17').Value = '  -2.43%  '
$ws.Range('D18').Value = "'434.70"
$ws.Range('E18').Value = '  -2.52%  '
$ws.Range('D19').Value = "'5.52"
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').Value = "'13.20"
$ws.Range('E20').Value = '  -3.24%  '
$ws.Range('D21').Value = "'7.40"
$ws.Range('E21').Value = '  -4.17%  '
$ws.Range('D22').Value = "'71.78"
$ws.Range('E22').Value = '  -2.98%  '
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '3.421.69'
$ws.Range('E24').Value = '  -1.95%  '
$ws.Range('D25').Value = "'0.506"
$ws.Range('E25').Value = '  -1.42%  '
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('D27').Value = "'0.0000113"
$ws.Range('E27').Value = '  -6.46%  '
$ws.Range('E28').Value = '  -2.37%  '
$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  -1.85%  '
$ws.Range('D31').Value = "'22.29"
$ws.Range('E31').Value = '  -3.07%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').Value = "'5.19"
$ws.Range('E33').Value = '  -3.20%  '
$ws.Range('E34').Value = '  -3.34%  '
$ws.Range('E35').Value = '  -4.41%  '
$ws.Range('D36').Value = "'156.81"
$ws.Range('E36').Value = '  -3.09%  '
$ws.Range('E37').Value = '  -4.97%  '
$ws.Range('D38').Value = "'26.54"
$ws.Range('E38').Value = '  -3.87%  '
$ws.Range('D39').Value = "'1.79"
$ws.Range('E39').Value = '  -3.58%  '
$ws.Range('D40').Value = '2.774.42'
$ws.Range('E40').Value = '  -2.08%  '
$ws.Range('D41').Value = "'0.775"
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('D42').Value = "'4.30"
$ws.Range('E42').Value = '  -4.23%  '
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').Value = "'6.06"
$ws.Range('E44').Value = '  -2.80%  '
$ws.Range('D45').Value = "'0.0659"
$ws.Range('E45').Value = '  -2.01%  '
$ws.Range('D46').Value = "'321.53"
$ws.Range('E46').Value = '  -0.90%  '
$ws.Range('D47').Value = "'2.29"
$ws.Range('E47').Value = '  -3.78%  '
$ws.Range('D48').Value = "'23.15"
$ws.Range('E48').Value = '  -5.99%  '
$ws.Range('D49').Value = "'0.0267"
$ws.Range('E49').Value = '  -2.79%  '
$ws.Range('E50').Value = '  +2.47%  '
$ws.Range('E51').Value = '  -0.02%  '
